$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing note text for B16/D16 row (shared string previously "Looking into confusion matrices to aid in precision")
$ws.Range("D16").Value = "Confusion matrix only useful for prediction, "

# B16 hours value changes from 5 to 6
$ws.Range("B16").Value = 6

# Add new row 17: date, hours, note
$ws.Range("A17").Value = 45430
$ws.Range("A17").NumberFormat = "d-mmm"
$ws.Range("B17").Value = 1
$ws.Range("D17").Value = "Adding in BA student to widen the model training pool"
$ws.Range("D17").WrapText = $true

# Update selection / view to match the new last row
$ws.Range("D17").Select()
